$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15: H15=2031.9572, I15=2031.9572, K15=6095.8716, M15=-5926.8716
$ws.Range("H15").Value = 2031.9572
$ws.Range("I15").Value = 2031.9572
$ws.Range("K15").Value = 6095.8716
$ws.Range("M15").Value = -5926.8716

# row 28: H28=1049.3636, I28=783, K28=783, M28=-298
$ws.Range("H28").Value = 1049.3636
$ws.Range("I28").Value = 783
$ws.Range("K28").Value = 783
$ws.Range("M28").Value = -298

# row 33: H33=2684.4167, I33=2369.889, J33=3628, K33=2369.889, L33=3628, M33=-2140.889, N33=-4086
$ws.Range("H33").Value = 2684.4167
$ws.Range("I33").Value = 2369.889
$ws.Range("J33").Value = 3628
$ws.Range("K33").Value = 2369.889
$ws.Range("L33").Value = 3628
$ws.Range("M33").Value = -2140.889
$ws.Range("N33").Value = -4086

# row 51: H51=21746.908, J51=23546.6, L51=23546.6, N51=-24514.6
$ws.Range("H51").Value = 21746.908
$ws.Range("J51").Value = 23546.6
$ws.Range("L51").Value = 23546.6
$ws.Range("N51").Value = -24514.6

# row 62: H62=4404.2856, I62=3594.625, K62=3594.625, M62=-2970.625
$ws.Range("H62").Value = 4404.2856
$ws.Range("I62").Value = 3594.625
$ws.Range("K62").Value = 3594.625
$ws.Range("M62").Value = -2970.625

# row 65: H65=4404.2856, I65=3594.625, K65=17973.125, M65=-14853.125
$ws.Range("H65").Value = 4404.2856
$ws.Range("I65").Value = 3594.625
$ws.Range("K65").Value = 17973.125
$ws.Range("M65").Value = -14853.125

# row 99: H99=513.6667, I99=496.81818, K99=1490.45454, M99=7.545460000000048
$ws.Range("H99").Value = 513.6667
$ws.Range("I99").Value = 496.81818
$ws.Range("K99").Value = 1490.45454
$ws.Range("M99").Value = 7.545460000000048

# row 113: H113=9333, I113=10000, J113=8999.5, K113=10000, L113=8999.5, M113=-6746, N113=-15507.5
$ws.Range("H113").Value = 9333
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 8999.5
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 8999.5
$ws.Range("M113").Value = -6746
$ws.Range("N113").Value = -15507.5

# row 137: H137=4672.76, I137=1378.125, K137=4134.375, M137=-1584.375
$ws.Range("H137").Value = 4672.76
$ws.Range("I137").Value = 1378.125
$ws.Range("K137").Value = 4134.375
$ws.Range("M137").Value = -1584.375

# row 138: H138=2761.9688, I138=1386.75, J138=3079.327, K138=4160.25, L138=9237.981, M138=979.75, N138=-19517.981
$ws.Range("H138").Value = 2761.9688
$ws.Range("I138").Value = 1386.75
$ws.Range("J138").Value = 3079.327
$ws.Range("K138").Value = 4160.25
$ws.Range("L138").Value = 9237.981
$ws.Range("M138").Value = 979.75
$ws.Range("N138").Value = -19517.981

$ws = $wb.Worksheets.Item("ARM")
# row 2: H2=1778.9, I2=1754.4445, K2=1754.4445, M2=-1641.4445
$ws.Range("H2").Value = 1778.9
$ws.Range("I2").Value = 1754.4445
$ws.Range("K2").Value = 1754.4445
$ws.Range("M2").Value = -1641.4445

# row 32: H32=7818971.5, I32=10206851, K32=10206851, M32=-10206564
$ws.Range("H32").Value = 7818971.5
$ws.Range("I32").Value = 10206851
$ws.Range("K32").Value = 10206851
$ws.Range("M32").Value = -10206564

# row 102: H102=5238.0415, I102=5708.7144, K102=5708.7144, M102=-4086.7144
$ws.Range("H102").Value = 5238.0415
$ws.Range("I102").Value = 5708.7144
$ws.Range("K102").Value = 5708.7144
$ws.Range("M102").Value = -4086.7144

# row 113: H113=66000, J113=66000, L113=66000, N113=-74678
$ws.Range("H113").Value = 66000
$ws.Range("J113").Value = 66000
$ws.Range("L113").Value = 66000
$ws.Range("N113").Value = -74678

# row 116: H116=1778.9, I116=1754.4445, K116=1754.4445, M116=539.5554999999999
$ws.Range("H116").Value = 1778.9
$ws.Range("I116").Value = 1754.4445
$ws.Range("K116").Value = 1754.4445
$ws.Range("M116").Value = 539.5554999999999

# row 132: H132=27112.762, I132=27488.975, J132=22222, K132=82466.92499999999, L132=66666, M132=-79936.92499999999, N132=-71726
$ws.Range("H132").Value = 27112.762
$ws.Range("I132").Value = 27488.975
$ws.Range("J132").Value = 22222
$ws.Range("K132").Value = 82466.92499999999
$ws.Range("L132").Value = 66666
$ws.Range("M132").Value = -79936.92499999999
$ws.Range("N132").Value = -71726

$ws = $wb.Worksheets.Item("BSM")
# row 3: H3=1778.9, I3=1754.4445, K3=1754.4445, M3=-1640.4445
$ws.Range("H3").Value = 1778.9
$ws.Range("I3").Value = 1754.4445
$ws.Range("K3").Value = 1754.4445
$ws.Range("M3").Value = -1640.4445

# row 20: H20=6429.6, I20=5935.5, J20=7417.8, K20=5935.5, L20=7417.8, M20=-5688.5, N20=-7911.8
$ws.Range("H20").Value = 6429.6
$ws.Range("I20").Value = 5935.5
$ws.Range("J20").Value = 7417.8
$ws.Range("K20").Value = 5935.5
$ws.Range("L20").Value = 7417.8
$ws.Range("M20").Value = -5688.5
$ws.Range("N20").Value = -7911.8

# row 82: H82=47881.332, J82=133000, L82=133000, N82=-133766
$ws.Range("H82").Value = 47881.332
$ws.Range("J82").Value = 133000
$ws.Range("L82").Value = 133000
$ws.Range("N82").Value = -133766

# row 85: H85=47881.332, J85=133000, L85=133000, N85=-135652
$ws.Range("H85").Value = 47881.332
$ws.Range("J85").Value = 133000
$ws.Range("L85").Value = 133000
$ws.Range("N85").Value = -135652

# row 86: H86=2684.5186, I86=2373.2104, J86=3423.875, K86=2373.2104, L86=3423.875, M86=-1250.2104, N86=-5669.875
$ws.Range("H86").Value = 2684.5186
$ws.Range("I86").Value = 2373.2104
$ws.Range("J86").Value = 3423.875
$ws.Range("K86").Value = 2373.2104
$ws.Range("L86").Value = 3423.875
$ws.Range("M86").Value = -1250.2104
$ws.Range("N86").Value = -5669.875

# row 89: H89=2684.5186, I89=2373.2104, J89=3423.875, K89=11866.052, L89=17119.375, M89=-6250.052, N89=-28351.375
$ws.Range("H89").Value = 2684.5186
$ws.Range("I89").Value = 2373.2104
$ws.Range("J89").Value = 3423.875
$ws.Range("K89").Value = 11866.052
$ws.Range("L89").Value = 17119.375
$ws.Range("M89").Value = -6250.052
$ws.Range("N89").Value = -28351.375

# row 107: H107=479.5, I107=434.2, K107=434.2, M107=1485.8
$ws.Range("H107").Value = 479.5
$ws.Range("I107").Value = 434.2
$ws.Range("K107").Value = 434.2
$ws.Range("M107").Value = 1485.8

$ws = $wb.Worksheets.Item("CRP")
# row 31: H31=659604.75, I31=12383.5, K31=12383.5, M31=-12088.5
$ws.Range("H31").Value = 659604.75
$ws.Range("I31").Value = 12383.5
$ws.Range("K31").Value = 12383.5
$ws.Range("M31").Value = -12088.5

# row 34: H34=659604.75, I34=12383.5, K34=12383.5, M34=-12181.5
$ws.Range("H34").Value = 659604.75
$ws.Range("I34").Value = 12383.5
$ws.Range("K34").Value = 12383.5
$ws.Range("M34").Value = -12181.5

# row 58: H58=7532.143, I58=5942.4, K58=5942.4, M58=-5739.4
$ws.Range("H58").Value = 7532.143
$ws.Range("I58").Value = 5942.4
$ws.Range("K58").Value = 5942.4
$ws.Range("M58").Value = -5739.4

# row 107: H107=3452.2, I107=1799, K107=1799, M107=121
$ws.Range("H107").Value = 3452.2
$ws.Range("I107").Value = 1799
$ws.Range("K107").Value = 1799
$ws.Range("M107").Value = 121

# row 132: H132=3832.0889, I132=1739.6571, J132=11155.6, K132=5218.971299999999, L132=33466.8, M132=-2688.971299999999, N132=-38526.8
$ws.Range("H132").Value = 3832.0889
$ws.Range("I132").Value = 1739.6571
$ws.Range("J132").Value = 11155.6
$ws.Range("K132").Value = 5218.971299999999
$ws.Range("L132").Value = 33466.8
$ws.Range("M132").Value = -2688.971299999999
$ws.Range("N132").Value = -38526.8

# row 133: H133=85000, J133=85000, L133=85000, N133=-90060
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -90060

# row 136: H136=7532.143, I136=5942.4, K136=17827.2, M136=-15277.2
$ws.Range("H136").Value = 7532.143
$ws.Range("I136").Value = 5942.4
$ws.Range("K136").Value = 17827.2
$ws.Range("M136").Value = -15277.2

$ws = $wb.Worksheets.Item("CUL")
# row 138: H138=3985.5386, I138=3116, K138=9348, M138=-4208
$ws.Range("H138").Value = 3985.5386
$ws.Range("I138").Value = 3116
$ws.Range("K138").Value = 9348
$ws.Range("M138").Value = -4208

# row 141: H141=225247, J141=14975, L141=44925, N141=-55285
$ws.Range("H141").Value = 225247
$ws.Range("J141").Value = 14975
$ws.Range("L141").Value = 44925
$ws.Range("N141").Value = -55285

$ws = $wb.Worksheets.Item("GSM")
# row 80: H80=4955.3687, J80=6448.5, L80=6448.5, N80=-8444.5
$ws.Range("H80").Value = 4955.3687
$ws.Range("J80").Value = 6448.5
$ws.Range("L80").Value = 6448.5
$ws.Range("N80").Value = -8444.5

# row 83: H83=4955.3687, J83=6448.5, L83=32242.5, N83=-42226.5
$ws.Range("H83").Value = 4955.3687
$ws.Range("J83").Value = 6448.5
$ws.Range("L83").Value = 32242.5
$ws.Range("N83").Value = -42226.5

# row 102: H102=1850.973, I102=956.6539, K102=956.6539, M102=665.3461
$ws.Range("H102").Value = 1850.973
$ws.Range("I102").Value = 956.6539
$ws.Range("K102").Value = 956.6539
$ws.Range("M102").Value = 665.3461

# row 113: H113=4169.684, J113=4502.8, L113=4502.8, N113=-8842.799999999999
$ws.Range("H113").Value = 4169.684
$ws.Range("J113").Value = 4502.8
$ws.Range("L113").Value = 4502.8
$ws.Range("N113").Value = -8842.799999999999

# row 122: H122=2955.8, I122=1819.75, K122=5459.25, M122=-3009.25
$ws.Range("H122").Value = 2955.8
$ws.Range("I122").Value = 1819.75
$ws.Range("K122").Value = 5459.25
$ws.Range("M122").Value = -3009.25

# row 126: H126=4043.2666, I126=3683.2222, J126=4583.3335, K126=11049.6666, L126=13750.0005, M126=-8579.6666, N126=-18690.0005
$ws.Range("H126").Value = 4043.2666
$ws.Range("I126").Value = 3683.2222
$ws.Range("J126").Value = 4583.3335
$ws.Range("K126").Value = 11049.6666
$ws.Range("L126").Value = 13750.0005
$ws.Range("M126").Value = -8579.6666
$ws.Range("N126").Value = -18690.0005

# row 132: H132=55558804, I132=66669970, K132=200009910, M132=-200007380
$ws.Range("H132").Value = 55558804
$ws.Range("I132").Value = 66669970
$ws.Range("K132").Value = 200009910
$ws.Range("M132").Value = -200007380

# row 136: H136=12379.2, J136=12379.2, L136=37137.60000000001, N136=-42237.60000000001
$ws.Range("H136").Value = 12379.2
$ws.Range("J136").Value = 12379.2
$ws.Range("L136").Value = 37137.60000000001
$ws.Range("N136").Value = -42237.60000000001

$ws = $wb.Worksheets.Item("LTW")
# row 40: H40=4751.5, I40=3061, J40=8857, K40=3061, L40=8857, M40=-2925, N40=-9129
$ws.Range("H40").Value = 4751.5
$ws.Range("I40").Value = 3061
$ws.Range("J40").Value = 8857
$ws.Range("K40").Value = 3061
$ws.Range("L40").Value = 8857
$ws.Range("M40").Value = -2925
$ws.Range("N40").Value = -9129

# row 61: H61=1274.7727, I61=984.2, K61=984.2, M61=-782.2
$ws.Range("H61").Value = 1274.7727
$ws.Range("I61").Value = 984.2
$ws.Range("K61").Value = 984.2
$ws.Range("M61").Value = -782.2

# row 113: H113=1274.7727, I113=984.2, K113=984.2, M113=1185.8
$ws.Range("H113").Value = 1274.7727
$ws.Range("I113").Value = 984.2
$ws.Range("K113").Value = 984.2
$ws.Range("M113").Value = 1185.8

# row 132: H132=299975.03, I132=6438.2334, J132=2501501, K132=19314.7002, L132=7504503, M132=-16784.7002, N132=-7509563
$ws.Range("H132").Value = 299975.03
$ws.Range("I132").Value = 6438.2334
$ws.Range("J132").Value = 2501501
$ws.Range("K132").Value = 19314.7002
$ws.Range("L132").Value = 7504503
$ws.Range("M132").Value = -16784.7002
$ws.Range("N132").Value = -7509563

# row 136: H136=52655.69, I136=7427.647, J136=138086.44, K136=22282.941, L136=414259.32, M136=-19732.941, N136=-419359.32
$ws.Range("H136").Value = 52655.69
$ws.Range("I136").Value = 7427.647
$ws.Range("J136").Value = 138086.44
$ws.Range("K136").Value = 22282.941
$ws.Range("L136").Value = 414259.32
$ws.Range("M136").Value = -19732.941
$ws.Range("N136").Value = -419359.32

$ws = $wb.Worksheets.Item("WVR")
# row 113: H113=1128, I113=999.3333, K113=2997.9999, M113=-827.9998999999998
$ws.Range("H113").Value = 1128
$ws.Range("I113").Value = 999.3333
$ws.Range("K113").Value = 2997.9999
$ws.Range("M113").Value = -827.9998999999998
